# Adjust scenario names in Power_BusInfo and Power_Network:
#   scenarioA -> ScenarioA
#   scenarioB -> ScenarioB
# Renaming the sheets automatically updates all formulas / defined names
# that reference them (_xlnm._FilterDatabase, "network", etc.).
# After renaming, make sure ScenarioA (the first sheet) is the active /
# selected tab, since ScenarioB was previously marked as selected.

$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("scenarioA")
$wsB = $wb.Worksheets.Item("scenarioB")

$wsA.Name = "ScenarioA"
$wsB.Name = "ScenarioB"

# ScenarioA is the first sheet and should end up as the selected tab
# (ScenarioB had previously been left marked as the selected tab).
$wsA.Activate()
